$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.691.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.419.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.88'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +8.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.419.45'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.43%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.477'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.126'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.79%  '
$ws.Range('E12').Value = '  +5.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.995.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.18%  '
$ws.Range('E14').Value = '  +2.02%  '
$ws.Range('E15').Value = '  +7.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.401.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.813.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('E19').Value = '  +5.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.88%  '
$ws.Range('E23').Value = '  +3.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.554.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.37%  '
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('E26').Value = '  +18.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.82'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.991'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.32'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.10%  '
$ws.Range('E32').Value = '  +6.51%  '
$ws.Range('E33').Value = '  +2.59%  '
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.71'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.53'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.22%  '
$ws.Range('E39').Value = '  +6.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.00'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0798'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +11.50%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.775'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.27%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.46'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.26%  '
$ws.Range('E49').Value = '  +5.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.91'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.357.13'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.37%  '
